# Work plan v. 2.0.
# Rename the worksheet from "Effort per Task" to "Effort by Task".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Effort by Task"
